$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$ws.Range("A1").Value = "Accounts"
$ws.Range("B1").Value = "1 Jul, 2023 - 30 Mar, 2024"
$ws.Range("C1").Value = "1 Jul, 2022 - 30 Mar, 2023 (PP)"
$ws.Range("D1").Value = "outliers"
$ws.Range("E1").Value = "Account type"
$ws.Range("F1").Value = "Account hierarchy"
$ws.Range("G1").Value = "Difference 1 Jul, 2023 - 30 Mar, 2024"
$ws.Range("H1").Value = "% of Sales 1 Jul, 2023 - 30 Mar, 2024"
$ws.Range("I1").Value = "% of Sales 1 Jul, 2022 - 30 Mar, 2023 (PP)"
$ws.Range("J1").Value = "% Difference 1 Jul, 2023 - 30 Mar, 2024"
$ws.Range("K1").Value = "Percentage of Sales Difference 1 Jul, 2023 - 30 Mar, 2024"

# --- Data rows (rows 2-21): full replacement with new Cost of Sales breakdown ---
$ws.Range("A2").Value = " Salaries- Operation Staff"
$ws.Range("B2").Value = 41346905
$ws.Range("C2").Value = 35983987.33
$ws.Range("D2").Value = $true
$ws.Range("E2").Value = "Cost of Sales"
$ws.Range("F2").Value = ""
$ws.Range("G2").Value = 5362917.670000002
$ws.Range("H2").Value = 12.4
$ws.Range("I2").Value = 10.8
$ws.Range("J2").Value = 14.9
$ws.Range("K2").Value = 1.6

$ws.Range("A3").Value = " Chicken"
$ws.Range("B3").Value = 31584301
$ws.Range("C3").Value = 36158003.28
$ws.Range("D3").Value = $true
$ws.Range("E3").Value = "Cost of Sales"
$ws.Range("F3").Value = ""
$ws.Range("G3").Value = -4573702.280000001
$ws.Range("H3").Value = 9.5
$ws.Range("I3").Value = 10.9
$ws.Range("J3").Value = -12.6
$ws.Range("K3").Value = -1.4

$ws.Range("A4").Value = " Mutton"
$ws.Range("B4").Value = 23775608
$ws.Range("C4").Value = 24046624
$ws.Range("D4").Value = $false
$ws.Range("E4").Value = "Cost of Sales"
$ws.Range("F4").Value = ""
$ws.Range("G4").Value = -271016
$ws.Range("H4").Value = 7.1
$ws.Range("I4").Value = 7.2
$ws.Range("J4").Value = -1.1
$ws.Range("K4").Value = -0.1

$ws.Range("A5").Value = " Utilitites"
$ws.Range("B5").Value = 17422851.79
$ws.Range("C5").Value = 12200143
$ws.Range("D5").Value = $true
$ws.Range("E5").Value = "Cost of Sales"
$ws.Range("F5").Value = ""
$ws.Range("G5").Value = 5222708.789999999
$ws.Range("H5").Value = 5.2
$ws.Range("I5").Value = 3.7
$ws.Range("J5").Value = 42.8
$ws.Range("K5").Value = 1.6

$ws.Range("A6").Value = " Other Dry Stock"
$ws.Range("B6").Value = 15561741
$ws.Range("C6").Value = 14325630.78
$ws.Range("D6").Value = $false
$ws.Range("E6").Value = "Cost of Sales"
$ws.Range("F6").Value = ""
$ws.Range("G6").Value = 1236110.220000001
$ws.Range("H6").Value = 4.7
$ws.Range("I6").Value = 4.3
$ws.Range("J6").Value = 8.6
$ws.Range("K6").Value = 0.4

$ws.Range("A7").Value = " Dairy & Bakery Items"
$ws.Range("B7").Value = 15304457
$ws.Range("C7").Value = 12683930.74
$ws.Range("D7").Value = $true
$ws.Range("E7").Value = "Cost of Sales"
$ws.Range("F7").Value = ""
$ws.Range("G7").Value = 2620526.26
$ws.Range("H7").Value = 4.6
$ws.Range("I7").Value = 3.8
$ws.Range("J7").Value = 20.7
$ws.Range("K7").Value = 0.8

$ws.Range("A8").Value = " Rents, rates & taxes"
$ws.Range("B8").Value = 13197737
$ws.Range("C8").Value = 10508556
$ws.Range("D8").Value = $true
$ws.Range("E8").Value = "Cost of Sales"
$ws.Range("F8").Value = ""
$ws.Range("G8").Value = 2689181
$ws.Range("H8").Value = 4
$ws.Range("I8").Value = 3.2
$ws.Range("J8").Value = 25.6
$ws.Range("K8").Value = 0.8

$ws.Range("A9").Value = " Fuel & Transporation"
$ws.Range("B9").Value = 12291882
$ws.Range("C9").Value = 13212231.3
$ws.Range("D9").Value = $false
$ws.Range("E9").Value = "Cost of Sales"
$ws.Range("F9").Value = ""
$ws.Range("G9").Value = -920349.3000000007
$ws.Range("H9").Value = 3.7
$ws.Range("I9").Value = 4
$ws.Range("J9").Value = -7
$ws.Range("K9").Value = -0.3

$ws.Range("A10").Value = " Vegetables & Fruits"
$ws.Range("B10").Value = 8669529
$ws.Range("C10").Value = 8412028
$ws.Range("D10").Value = $false
$ws.Range("E10").Value = "Cost of Sales"
$ws.Range("F10").Value = ""
$ws.Range("G10").Value = 257501
$ws.Range("H10").Value = 2.6
$ws.Range("I10").Value = 2.5
$ws.Range("J10").Value = 3.1
$ws.Range("K10").Value = 0.1

$ws.Range("A11").Value = " Oil"
$ws.Range("B11").Value = 7376256
$ws.Range("C11").Value = 8693172.76
$ws.Range("D11").Value = $false
$ws.Range("E11").Value = "Cost of Sales"
$ws.Range("F11").Value = ""
$ws.Range("G11").Value = -1316916.76
$ws.Range("H11").Value = 2.2
$ws.Range("I11").Value = 2.6
$ws.Range("J11").Value = -15.1
$ws.Range("K11").Value = -0.4

$ws.Range("A12").Value = " Wages & Allowance"
$ws.Range("B12").Value = 7100965
$ws.Range("C12").Value = 8354416
$ws.Range("D12").Value = $false
$ws.Range("E12").Value = "Cost of Sales"
$ws.Range("F12").Value = ""
$ws.Range("G12").Value = -1253451
$ws.Range("H12").Value = 2.1
$ws.Range("I12").Value = 2.5
$ws.Range("J12").Value = -15
$ws.Range("K12").Value = -0.4

$ws.Range("A13").Value = " Outsourced Food"
$ws.Range("B13").Value = 6995449
$ws.Range("C13").Value = 9884767
$ws.Range("D13").Value = $true
$ws.Range("E13").Value = "Cost of Sales"
$ws.Range("F13").Value = ""
$ws.Range("G13").Value = -2889318
$ws.Range("H13").Value = 2.1
$ws.Range("I13").Value = 3
$ws.Range("J13").Value = -29.2
$ws.Range("K13").Value = -0.9

$ws.Range("A14").Value = " Labour Out Sourced"
$ws.Range("B14").Value = 6910950
$ws.Range("C14").Value = 5625265
$ws.Range("D14").Value = $false
$ws.Range("E14").Value = "Cost of Sales"
$ws.Range("F14").Value = ""
$ws.Range("G14").Value = 1285685
$ws.Range("H14").Value = 2.1
$ws.Range("I14").Value = 1.7
$ws.Range("J14").Value = 22.9
$ws.Range("K14").Value = 0.4

$ws.Range("A15").Value = " Other Consumables"
$ws.Range("B15").Value = 6536746.6
$ws.Range("C15").Value = 7117116
$ws.Range("D15").Value = $false
$ws.Range("E15").Value = "Cost of Sales"
$ws.Range("F15").Value = ""
$ws.Range("G15").Value = -580369.4000000004
$ws.Range("H15").Value = 2
$ws.Range("I15").Value = 2.1
$ws.Range("J15").Value = -8.2
$ws.Range("K15").Value = -0.2

$ws.Range("A16").Value = " Rice"
$ws.Range("B16").Value = 5914882
$ws.Range("C16").Value = 5769558
$ws.Range("D16").Value = $false
$ws.Range("E16").Value = "Cost of Sales"
$ws.Range("F16").Value = ""
$ws.Range("G16").Value = 145324
$ws.Range("H16").Value = 1.8
$ws.Range("I16").Value = 1.7
$ws.Range("J16").Value = 2.5
$ws.Range("K16").Value = 0

$ws.Range("A17").Value = " Beef"
$ws.Range("B17").Value = 5784255
$ws.Range("C17").Value = 4459168
$ws.Range("D17").Value = $false
$ws.Range("E17").Value = "Cost of Sales"
$ws.Range("F17").Value = ""
$ws.Range("G17").Value = 1325087
$ws.Range("H17").Value = 1.7
$ws.Range("I17").Value = 1.3
$ws.Range("J17").Value = 29.7
$ws.Range("K17").Value = 0.4

$ws.Range("A18").Value = " Disposible Material"
$ws.Range("B18").Value = 5742937
$ws.Range("C18").Value = 7373230
$ws.Range("D18").Value = $false
$ws.Range("E18").Value = "Cost of Sales"
$ws.Range("F18").Value = ""
$ws.Range("G18").Value = -1630293
$ws.Range("H18").Value = 1.7
$ws.Range("I18").Value = 2.2
$ws.Range("J18").Value = -22.1
$ws.Range("K18").Value = -0.5

$ws.Range("A19").Value = " Beverages & Soft drinks"
$ws.Range("B19").Value = 5594702
$ws.Range("C19").Value = 4556449.65
$ws.Range("D19").Value = $false
$ws.Range("E19").Value = "Cost of Sales"
$ws.Range("F19").Value = ""
$ws.Range("G19").Value = 1038252.35
$ws.Range("H19").Value = 1.7
$ws.Range("I19").Value = 1.4
$ws.Range("J19").Value = 22.8
$ws.Range("K19").Value = 0.3

$ws.Range("A20").Value = " Fish & Prawns"
$ws.Range("B20").Value = 4083157
$ws.Range("C20").Value = 4295326
$ws.Range("D20").Value = $false
$ws.Range("E20").Value = "Cost of Sales"
$ws.Range("F20").Value = ""
$ws.Range("G20").Value = -212169
$ws.Range("H20").Value = 1.2
$ws.Range("I20").Value = 1.3
$ws.Range("J20").Value = -4.9
$ws.Range("K20").Value = -0.1

$ws.Range("A21").Value = " Flour"
$ws.Range("B21").Value = 4023287
$ws.Range("C21").Value = 3890931
$ws.Range("D21").Value = $false
$ws.Range("E21").Value = "Cost of Sales"
$ws.Range("F21").Value = ""
$ws.Range("G21").Value = 132356
$ws.Range("H21").Value = 1.2
$ws.Range("I21").Value = 1.2
$ws.Range("J21").Value = 3.4
$ws.Range("K21").Value = 0

